# Applies the "generate by ai feature" revision to the RTOS deck:
#  - Slide 2 title/subtitle wording tweak
#  - Slide 3 (Index) bullet list rewritten for the new outline
#  - Slides 4-9 retitled/rebodied to match the revised topic flow
#  - Slide 10 repurposed from "Conclusion" into "Choosing the Right RTOS
#    and Development Tools"
#  - A brand-new Slide 11 added at the end holding the real conclusion
#    (cloned from the original Slide 10 so its paragraph/run structure
#    matches, then retitled/rebodied)

$p = $ppt.ActivePresentation

# --- Slide 11 (new): duplicate the original "Conclusion" slide (slide 10)
#     BEFORE editing slide 10's text, so the clone keeps the exact same
#     txBody shape (leading empty paragraph, <a:pPr/>, single run) and
#     lands at the end of the deck.
$s10_orig = $p.Slides.Item(10)
$dupRange = $s10_orig.Duplicate()
$s11 = $p.Slides.Item(11)
[void]$s11.Shapes.Item(1).TextFrame.TextRange.Replace("Conclusion: Summary and Future Trends", "Conclusion and Key Takeaways")
[void]$s11.Shapes.Item(2).TextFrame.TextRange.Replace("Summary of Key Concepts: Recap of RTOS features, architectures, and common implementations.Takeaways:  Understanding the importance of RTOS in real-time applications.Future Trends:  Evolution of RTOS towards AI integration, cloud connectivity, and enhanced security.Q&A: Open the floor for questions and discussion.", "Summary of key concepts covered in the presentation.Importance of RTOS in modern embedded systems.Recap of the advantages and challenges of using RTOS.Future trends in RTOS technology.Call to action: Encourage further exploration and learning.Q&A session.")

# --- Slide 2: title slide -------------------------------------------------
$s2 = $p.Slides.Item(2)
[void]$s2.Shapes.Item(1).TextFrame.TextRange.Replace("Unveiling the Power of RTOS: Real-Time Operating Systems", "Unlocking the Power of RTOS: Real-Time Operating Systems")
[void]$s2.Shapes.Item(2).TextFrame.TextRange.Replace("Created by: Aditya Bhogil", "Created by: mayur")

# --- Slide 3: index / table of contents ------------------------------------
$s3 = $p.Slides.Item(3)
[void]$s3.Shapes.Item(2).TextFrame.TextRange.Replace("1.Introduction to Real-Time Operating Systems (RTOS)2.Key Features and Benefits of RTOS3.RTOS Architecture: Tasks, Scheduling, and Inter-process Communication4.Common RTOS Implementations and Use Cases5.Choosing the Right RTOS for Your Application6.RTOS Development and Debugging Techniques7.Conclusion: Summary and Future TrendsConclusion", "1.Introduction to RTOS2.Key Concepts and Terminology3.RTOS Scheduling Algorithms4.Memory Management in RTOS5.Inter-Process Communication (IPC) Mechanisms6.RTOS Applications and Use Cases7.Choosing the Right RTOS and Development Tools8.Conclusion and Key TakeawaysConclusion")

# --- Slide 4: Introduction --------------------------------------------------
$s4 = $p.Slides.Item(4)
[void]$s4.Shapes.Item(1).TextFrame.TextRange.Replace("Introduction to Real-Time Operating Systems (RTOS)", "Introduction to RTOS")
[void]$s4.Shapes.Item(2).TextFrame.TextRange.Replace("Definition of an RTOS: A multitasking OS designed for real-time applications.Distinguishing RTOS from General-Purpose OS: Emphasis on deterministic behavior and timely response.Importance of Real-Time Constraints: Meeting deadlines is critical.Applications of RTOS: Automotive, Industrial Automation, Aerospace, Medical Devices.Key Characteristics: Predictability, Reliability, Efficiency.", "Definition and purpose of Real-Time Operating Systems (RTOS).Distinguishing characteristics of RTOS vs. general-purpose operating systems.Brief history and evolution of RTOS.Importance of determinism and predictability in RTOS.Advantages of using an RTOS in embedded systems.Examples of popular RTOS platforms (e.g., FreeRTOS, VxWorks, QNX).High-level overview of RTOS architecture.")

# --- Slide 5: Key Concepts and Terminology ---------------------------------
$s5 = $p.Slides.Item(5)
[void]$s5.Shapes.Item(1).TextFrame.TextRange.Replace("Key Features and Benefits of RTOS", "Key Concepts and Terminology")
[void]$s5.Shapes.Item(2).TextFrame.TextRange.Replace("Multitasking: Concurrent execution of multiple tasks.Preemptive Scheduling: Prioritized task execution.Inter-Process Communication (IPC): Mechanisms for task synchronization and data exchange (e.g., semaphores, mutexes, ...Memory Management: Efficient allocation and deallocation of memory resources.Real-Time Capabilities: Guaranteed response times within specified deadlines.Deterministic Behavior: Predictable system behavior.Improved System Reliability: Robust error handling and fault tolerance.", "Tasks/Threads: Defining and managing concurrent processes.Processes vs. Threads: Understanding the differences.Scheduling:  Prioritization and task execution.Context Switching: Efficient task transitions.Preemption: Interrupt handling and task interruption.Semaphores and Mutexes: Synchronization mechanisms.Deadlocks and race conditions: Understanding and preventing these issues.")

# --- Slide 6: RTOS Scheduling Algorithms -----------------------------------
$s6 = $p.Slides.Item(6)
[void]$s6.Shapes.Item(1).TextFrame.TextRange.Replace("RTOS Architecture: Tasks, Scheduling, and Inter-process Communication", "RTOS Scheduling Algorithms")
[void]$s6.Shapes.Item(2).TextFrame.TextRange.Replace("Task Management: Creation, deletion, and management of tasks.Scheduling Algorithms: Round-robin, priority-based, rate-monotonic.Context Switching: Efficient switching between tasks.Inter-Process Communication (IPC) Mechanisms: Semaphores, mutexes, message queues, mailboxes.Synchronization Techniques: Avoiding race conditions and deadlocks.Interrupt Handling: Efficient handling of hardware interrupts.", "Round Robin Scheduling: Fair but potentially inefficient.Priority-Based Scheduling: Efficient but requires careful priority assignment.Rate Monotonic Scheduling (RMS):  Real-time scheduling algorithm.Earliest Deadline First (EDF): Another real-time scheduling algorithm.Comparison of different scheduling algorithms: Pros and cons of each.Impact of scheduling algorithms on system performance and determinism.Choosing the appropriate scheduling algorithm for specific applications.")

# --- Slide 7: Memory Management in RTOS ------------------------------------
$s7 = $p.Slides.Item(7)
[void]$s7.Shapes.Item(1).TextFrame.TextRange.Replace("Common RTOS Implementations and Use Cases", "Memory Management in RTOS")
[void]$s7.Shapes.Item(2).TextFrame.TextRange.Replace("FreeRTOS: A widely used, open-source RTOS.VxWorks: A commercial RTOS known for its reliability and performance.QNX: A robust RTOS used in safety-critical applications.ThreadX: A royalty-free RTOS with a small footprint.Embedded Linux: A modified Linux kernel for embedded systems (often not strictly an RTOS).Use Case Examples:  Automotive engine control, industrial robotics, medical imaging equipment, flight control systems.", "Memory allocation and deallocation strategies.Static vs. dynamic memory allocation.Memory fragmentation and its impact.Memory protection mechanisms.Heap management techniques.Memory partitioning and its benefits.Techniques for optimizing memory usage in resource-constrained environments.")

# --- Slide 8: Inter-Process Communication (IPC) Mechanisms -----------------
$s8 = $p.Slides.Item(8)
[void]$s8.Shapes.Item(1).TextFrame.TextRange.Replace("Choosing the Right RTOS for Your Application", "Inter-Process Communication (IPC) Mechanisms")
[void]$s8.Shapes.Item(2).TextFrame.TextRange.Replace("Requirements Analysis: Defining real-time constraints, memory limitations, and processing power.RTOS Feature Comparison: Evaluating features such as scheduling algorithms, IPC mechanisms, and memory management.Scalability and Performance: Considering the future growth and performance needs of the application.Cost and Licensing: Evaluating the cost of the RTOS and its licensing model.Support and Community: Assessing the availability of support and community resources.Certification and Standards Compliance: Ensuring the RTOS meets relevant safety and security standards.", "Importance of IPC in multi-tasking environments.Message queues: Asynchronous communication.Semaphores: Synchronization and mutual exclusion.Mailboxes:  Structured data exchange.Shared memory: Efficient but requires careful synchronization.Pipes and sockets: Communication between processes.Choosing the appropriate IPC mechanism for specific needs.")

# --- Slide 9: RTOS Applications and Use Cases ------------------------------
$s9 = $p.Slides.Item(9)
[void]$s9.Shapes.Item(1).TextFrame.TextRange.Replace("RTOS Development and Debugging Techniques", "RTOS Applications and Use Cases")
[void]$s9.Shapes.Item(2).TextFrame.TextRange.Replace("Development Tools and IDEs: Selecting appropriate tools for development and debugging.Debugging Strategies: Using real-time debuggers and tracing tools.Real-Time Analysis: Analyzing system performance and identifying bottlenecks.Testing and Verification: Ensuring the RTOS meets the required real-time constraints and reliability requirements.Profiling and Optimization: Optimizing the RTOS and application code for performance.", "Automotive systems (e.g., Engine control, ABS).Industrial automation (e.g., Robotics, PLC).Medical devices (e.g., Pacemakers, imaging equipment).Aerospace and defense (e.g., Flight control, navigation).Consumer electronics (e.g., Smartwatches, smartphones).Real-time data acquisition and control systems.Networked embedded systems.")

# --- Slide 10: repurposed into "Choosing the Right RTOS and Development Tools"
$s10 = $p.Slides.Item(10)
[void]$s10.Shapes.Item(1).TextFrame.TextRange.Replace("Conclusion: Summary and Future Trends", "Choosing the Right RTOS and Development Tools")
[void]$s10.Shapes.Item(2).TextFrame.TextRange.Replace("Summary of Key Concepts: Recap of RTOS features, architectures, and common implementations.Takeaways:  Understanding the importance of RTOS in real-time applications.Future Trends:  Evolution of RTOS towards AI integration, cloud connectivity, and enhanced security.Q&A: Open the floor for questions and discussion.", "Factors to consider when selecting an RTOS (e.g., Real-time requirements, memory constraints, cost, support).Evaluating different RTOS platforms based on features and capabilities.Overview of popular RTOS development tools and IDEs.Importance of debugging and testing in RTOS development.Resources for learning and using RTOS.Open-source vs. commercial RTOS options.")
